# Add an "EnemyId" column (E) to the SceneQuest config table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the existing table to include the new column first, so Excel
# picks up the header text typed below as the new column's name.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:E4"))

# Header / meta rows for the new column
$ws.Range("E1").Value = "EnemyId"
$ws.Range("E2").Value = "int"
$ws.Range("E3").Value = "怪物id"
$ws.Range("E4").Value = 43000001

# Copy styles from existing row2/row3 cells into the new column cells
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection like in the target workbook
$ws.Range("E4").Select() | Out-Null
